$wb = $excel.ActiveWorkbook

# --- WebTables sheet: add new column E with data ---
$wsWebTables = $wb.Worksheets.Item("WebTables")

$wsWebTables.Range("E1").Value = "1nvalid"
$wsWebTables.Range("E2").Value = "1nvalid"
$wsWebTables.Range("E3").Value = "1nvalid"
$wsWebTables.Range("E4").Value = "aa"
$wsWebTables.Range("E5").Value = "aa"
$wsWebTables.Range("E6").Value = "aa"
$wsWebTables.Range("E7").Value = "1nvalid"

# --- TextBox sheet: move selection, no longer the active tab ---
$wsTextBox = $wb.Worksheets.Item("TextBox")
$wsTextBox.Activate()
$wsTextBox.Range("C30").Select()

# --- WebTables sheet becomes the active tab, with a new selection ---
$wsWebTables.Activate()
$wsWebTables.Range("F26").Select()
